$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: LinearRegression - only B2 changes
$ws.Range("B2").Value = 0.52388209438133

# Row 3: RandomForestRegressor - B3, C3, D3 change
$ws.Range("B3").Value = 0.03048732093053691
$ws.Range("C3").Value = 0.03051507313546532
$ws.Range("D3").Value = 0.06413505045966662

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor, values change
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.02970979745812575
$ws.Range("C4").Value = 0.02982992260983418
$ws.Range("D4").Value = 0.05749184399892781

# Row 5: AdaBoostRegressor -> MLPRegressor, values change
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.02538666422943526
$ws.Range("C5").Value = 0.02861434234868026
$ws.Range("D5").Value = 0.03366883031565795
